# chgt xxx colonne export + selection ihm
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Profil" column header (cell B1) to "ProfilXXXXX"
$ws.Range("B1").Value = "ProfilXXXXX"

# Update the remembered selection / active cell on the sheet (was D14, now D12)
$null = $ws.Range("D12").Select()
